# Clear form-response values in row 2 (all price/answer columns) while keeping
# columns A-D (timestamp/date/name/market) intact. Mirrors the author's
# "Add files via upload" edit that wiped the sample numeric answers from the
# second response row before re-uploading the form template.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToClear = @(
    "E2","F2","G2","H2","I2","J2","K2","L2","N2","O2",
    "Q2","T2","U2","W2","X2","Z2","AC2","AD2","AF2","AG2",
    "AI2","AJ2","AL2","AM2","AO2","AP2","AQ2","AR2","AS2","AT2",
    "AU2","AV2","AX2","AY2","AZ2","BA2","BB2","BC2","BD2","BG2",
    "BH2","BI2","BJ2","BK2","BL2","BM2","BP2","BQ2","BR2","BS2",
    "BV2","BW2","BY2","CB2","CC2","CD2","CE2","CH2","CI2","CJ2",
    "CK2","CN2","CO2","CP2","CQ2","CR2","CS2","CT2","CU2","CV2",
    "CW2","CX2","DA2","DB2","DC2","DD2","DE2","DF2","DG2","DM2",
    "DN2","DO2","DP2","DQ2","DS2","DV2","DW2","DY2","DZ2","EA2",
    "EB2","EE2","EF2","EH2","EI2","EJ2","EK2","EL2","EM2","EN2",
    "EO2","EP2","EQ2","ER2","ES2","ET2","EU2","EV2","EW2","EY2",
    "EZ2","FB2","FC2","FD2","FE2","FF2","FG2","FH2","FI2","FJ2",
    "FK2","FL2","FM2","FN2","FO2","FP2","FQ2","FR2","FT2"

)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

# Restore the cursor/selection to D2, matching the saved view state.
$ws.Range("D2").Select()
